$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($row = 2; $row -le 35; $row++) {
        $cell = $ws.Cells.Item($row, 4)
        if ($cell.Value2 -eq "(0, 0)") {
            $cell.Value2 = "(nan, nan)"
        }
    }
}
